# Add speaker notes ("commentaar") to the presentation, as described by the
# commit "update ppt (extra info in commentaar)".
#
# Slide 2 (the "Controller" remote-app screenshot) gets a multi-paragraph
# explanation of the textbox/logging, HDMI/TV sources and the two separate
# apps used for TV + remote control.
#
# Slide 3 (the "MainWindow" TV-simulation screenshot) gets an (empty) notes
# page created, with no additional text.

$p = $ppt.ActivePresentation

# --- Slide 2: detailed speaker notes -------------------------------------
$slide2 = $p.Slides.Item(2)
$notes2 = $slide2.NotesPage

$notesText = "Textbox: logging. `n" +
    "Source: 2x hdmi + TV`n" +
    "Settings: standaard zaken instellen. (Volume, channel, source bij opstarten tv)`n" +
    "`n" +
    "2 afzonderlijke applicaties. (2 solutions 1 TV en 1 afstandsbediening)."

$notes2.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = $notesText

# --- Slide 3: ensure a (blank) notes page exists --------------------------
$slide3 = $p.Slides.Item(3)
$notes3 = $slide3.NotesPage
$notes3.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = ""
